# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.085.13"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.650.86"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.44"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5208"
$ws.Range("E6").Value = "  -2.48%  "
$ws.Range("E7").Value = "  -0.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2616"
$ws.Range("E8").Value = "  -1.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06291"
$ws.Range("E9").Value = "  -1.93%  "
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07793"
$ws.Range("E11").Value = "  -0.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.477"
$ws.Range("E12").Value = "  -2.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.589.22"
$ws.Range("E13").Value = "  -4.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.877.99"
$ws.Range("E14").Value = "  -0.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5541"
$ws.Range("E15").Value = "  +0.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8003"
$ws.Range("E16").Value = "  -2.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.80"
$ws.Range("E17").Value = "  -1.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.079.20"
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.624"
$ws.Range("E20").Value = "  -1.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "194.15"
$ws.Range("E22").Value = "  -1.50%  "
$ws.Range("E23").Value = "  -1.59%  "
$ws.Range("E24").Value = "  -0.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.56"
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("E26").Value = "  -2.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.170"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.87"
$ws.Range("E28").Value = "  -1.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.477"
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05598"
$ws.Range("E30").Value = "  -4.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.475"
$ws.Range("E32").Value = "  -4.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.344"
$ws.Range("E33").Value = "  +1.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.595"
$ws.Range("E34").Value = "  -0.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.799"
$ws.Range("E35").Value = "  -0.94%  "
$ws.Range("E36").Value = "  -1.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.409"
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("E38").Value = "  -2.67%  "
$ws.Range("E39").Value = "  -1.56%  "
$ws.Range("E40").Value = "  +1.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.056.75"
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.006"
$ws.Range("E42").Value = "  -0.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8411"
$ws.Range("E43").Value = "  -2.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.53"
$ws.Range("E44").Value = "  -1.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.789.48"
$ws.Range("E45").Value = "  -0.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.18"
$ws.Range("E46").Value = "  -0.92%  "
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.008"
$ws.Range("E47").Value = "  -0.53%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₈104"
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05304"
$ws.Range("E49").Value = "  +2.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4336"
$ws.Range("E50").Value = "  -1.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.943"
$ws.Range("E51").Value = "  -1.49%  "
